$wb = $excel.ActiveWorkbook

# --- Sheet 1: testdata ---
# Column A (name) is edited top-to-bottom first, then column B (email)
# top-to-bottom, matching the order the shared strings were appended in.
$ws1 = $wb.Worksheets.Item("testdata")

$ws1.Range("A9").Value = "testEngineer4"
$ws1.Range("A10").Value = "testEngineer5"
$ws1.Range("A11").Value = "testEngineer6"
$ws1.Range("B9").Value = "test4@mailinator.com"
$ws1.Range("B10").Value = "test5@mailinator.com"
$ws1.Range("B11").Value = "test6@mailinator.com"

# --- Sheet 2: cxCreationValidKey ---
$ws2 = $wb.Worksheets.Item("cxCreationValidKey")

$ws2.Range("A2").Value = "testEngineer7"
$ws2.Range("A3").Value = "testEngineer8"
$ws2.Range("B2").Value = "test7@mailinator.com"
$ws2.Range("B3").Value = "test8@mailinator.com"

# --- Sheet 3: cxCreationInvalidKey ---
$ws3 = $wb.Worksheets.Item("cxCreationInvalidKey")

$ws3.Range("A2").Value = "testEngineer9"
$ws3.Range("A3").Value = "testEngineer10"
$ws3.Range("B2").Value = "test9@mailinator.com"
$ws3.Range("B3").Value = "test10@mailinator.com"

# --- Selections / active sheet / view state ---
# sheet1: drop the frozen "topLeftCell=A6" scroll anchor and move the
# cell-selection to B11
$ws1.Range("B11").Select()

# sheet2: selection moves to B3 (it also loses the tabSelected flag because
# sheet3 becomes the active tab below)
$ws2.Range("B3").Select()

# sheet3: becomes the active / selected tab, with C10:C11 selected
$ws3.Activate()
$ws3.Range("C10:C11").Select()
